# Settings sheet rework:
#  - remove the OutputFilePath / OutputSheetName rows
#  - add a new MailRecipient row (with a mailto: hyperlink) at the bottom
#  - re-establish the two pre-existing hyperlinks against their shifted cells
#  - select B11 (matches what was recorded when the workbook was saved)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Drop the OutputFilePath / OutputSheetName rows (old rows 3 & 4).
#    Everything below shifts up two rows:
#      old row 9  (YahooFinanceUrl)      -> new row 7
#      old row 10 (ExtractDataFromUrl)   -> new row 8
#      old row 11 (ExtractLimitFromURL)  -> new row 9
$ws.Rows("3:4").Delete()

# remember the real cell text so we can restore it after fixing up the
# hyperlink display text below (see step 2)
$extractUrlValue = $ws.Range("B8").Text

# 2) The hyperlink anchors are not shifted automatically when rows move, so
#    the worksheet's Hyperlinks collection still thinks the two links live on
#    the old row numbers. Clear everything out and rebuild against the
#    correct (post-shift) cells.
$ws.Range("A1").Hyperlinks.Delete()

# ExtractDataFromUrl row -> link to the "most active" listing (this mirrors
# the original file, where the hyperlink's display text differs from the
# cell's own JSON text). Added first so it keeps rId1, matching the source.
$linkB8 = $ws.Hyperlinks.Add($ws.Range("B8"), "https://ca.finance.yahoo.com/markets/stocks/most-active/")
$linkB8.TextToDisplay = "https://ca.finance.yahoo.com/markets/stocks/most-active/"
$ws.Range("B8").Value = $extractUrlValue
$ws.Range("B8").Style = "Hyperlink"

# YahooFinanceUrl row -> link back to the Yahoo Finance homepage (rId2)
$ws.Hyperlinks.Add($ws.Range("B7"), "https://finance.yahoo.com/")
$ws.Range("B7").Style = "Hyperlink"

# 3) Add the new MailRecipient row (row 10) - value first, then label, so the
#    shared-strings table ends up in the same order as the source workbook.
$ws.Range("B10").Value = "matan10cohen@gmail.com"
$ws.Hyperlinks.Add($ws.Range("B10"), "mailto:matan10cohen@gmail.com")
$ws.Range("B10").Style = "Hyperlink"

$ws.Range("A10").Value = "MailRecipient"
$ws.Range("A10").Font.Name = "Arial Unicode MS"
$ws.Range("A10").Font.Size = 10

# 4) Selection, as recorded in the saved workbook.
$ws.Range("B11").Select()
